# New results update (2021/03/27 12:16)
# Row 14's score (cell B15) was recomputed/corrected, which in turn
# changes the AVERAGE() result in B47. The cell also now carries Excel's
# built-in "Neutral" cell style (yellow) instead of "Bad" (red), matching
# the updated, less-poor score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Excel's built-in "Neutral" cell style (适中) to B15, reusing the
# workbook's existing style entry, then write the corrected value.
$ws.Range("B15").Style = "适中"
$ws.Range("B15").Value = 0.82730000000000004

# Refresh the view: drop the old scroll/zoom state and multi-cell
# selection, landing on C1 at 100% zoom like the freshly-reopened sheet.
$excel.ActiveWindow.Zoom = 100
$ws.Range("C1").Select() | Out-Null
